$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at BV (74th column). This shifts the existing
# "nom" (BV) and "url_produit" (BW) columns one place to the right,
# landing on BW and BX respectively, and grows the used range from
# BW206 to BX206 automatically.
$ws.Columns("BV").Insert()

# New column header: the latest scrape timestamp.
$ws.Range("BV1").Value = "2026-01-31 00:56:13"

# For every product row that still had a price in the last existing
# timestamp column (BU, rows 2-80), carry that same price forward into
# the freshly inserted BV column - i.e. record "price unchanged at this
# scrape time". Rows 81-206 have no price history left (BU is blank for
# them), so their new BV cell is left blank too.
for ($r = 2; $r -le 80; $r++) {
    $lastPrice = $ws.Cells.Item($r, 73).Value()
    $ws.Cells.Item($r, 74).Value = $lastPrice
}
